# Newton-Raphson change matrix workbook: add "OnePvAndOnePqNode" test sheet,
# and clean up the A22/A23 cells + selection on "TwoPqAndOnePvNodeDifferentOrde2".

$wb = $excel.ActiveWorkbook

# --- 1. Tidy up the existing "TwoPqAndOnePvNodeDifferentOrde2" sheet -------
$ws4 = $wb.Worksheets.Item("TwoPqAndOnePvNodeDifferentOrde2")
$ws4.Activate()
$ws4.Range("A22").ClearContents()
$ws4.Range("A23").ClearContents()
$ws4.Range("A22:A23").Select()

# --- 2. Add the new "OnePvAndOnePqNode" worksheet as the last tab ----------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $last)
$ws.Name = "OnePvAndOnePqNode"

# Block 1: generic admittances example (rows 1-3)
$ws.Range("A1").Value = "admittances"
$ws.Range("D1").Value = "real"
$ws.Range("G1").Value = "imaginary"
$ws.Range("J1").Value = "magnitude"
$ws.Range("M1").Value = "phase"

$ws.Range("A2").Formula = "=COMPLEX(0.6, -0.6)"
$ws.Range("B2").Formula = "=COMPLEX(-0.2, 0.4)"
$ws.Range("D2").Formula = "=IMREAL(A2)"
$ws.Range("E2").Formula = "=IMREAL(B2)"
$ws.Range("G2").Formula = "=IMAGINARY(A2)"
$ws.Range("H2").Formula = "=IMAGINARY(B2)"
$ws.Range("J2").Formula = "=IMABS(A2)"
$ws.Range("K2").Formula = "=IMABS(B2)"
$ws.Range("M2").Formula = "=IMARGUMENT(A2)"
$ws.Range("N2").Formula = "=IMARGUMENT(B2)"

$ws.Range("A3").Formula = "=COMPLEX(-0.2,0.4)"
$ws.Range("B3").Formula = "=COMPLEX(0.2,-0.4)"
$ws.Range("D3").Formula = "=IMREAL(A3)"
$ws.Range("E3").Formula = "=IMREAL(B3)"
$ws.Range("G3").Formula = "=IMAGINARY(A3)"
$ws.Range("H3").Formula = "=IMAGINARY(B3)"
$ws.Range("J3").Formula = "=IMABS(A3)"
$ws.Range("K3").Formula = "=IMABS(B3)"
$ws.Range("M3").Formula = "=IMARGUMENT(A3)"
$ws.Range("N3").Formula = "=IMARGUMENT(B3)"

# Block 2: voltages (rows 5-7)
$ws.Range("A5").Value = "voltages"
$ws.Range("C5").Value = "real"
$ws.Range("E5").Value = "imaginary"
$ws.Range("G5").Value = "magnitude"
$ws.Range("I5").Value = "phase"

$ws.Range("A6").Formula = "=COMPLEX(10,0)"
$ws.Range("C6").Formula = "=IMREAL(A6)"
$ws.Range("E6").Formula = "=IMAGINARY(A6)"
$ws.Range("G6").Formula = "=IMABS(A6)"
$ws.Range("I6").Formula = "=IMARGUMENT(A6)"

$ws.Range("A7").Formula = "=COMPLEX(10,0)"
$ws.Range("C7").Formula = "=IMREAL(A7)"
$ws.Range("E7").Formula = "=IMAGINARY(A7)"
$ws.Range("G7").Formula = "=IMABS(A7)"
$ws.Range("I7").Formula = "=IMARGUMENT(A7)"

# Block 3: currents (rows 9-11)
$ws.Range("A9").Value = "currents"
$ws.Range("C9").Value = "real"
$ws.Range("E9").Value = "imaginary"
$ws.Range("G9").Value = "magnitude"
$ws.Range("I9").Value = "phase"

$ws.Range("A10").Formula = "=COMPLEX(4.02,-1.96)"
$ws.Range("C10").Formula = "=IMREAL(A10)"
$ws.Range("E10").Formula = "=IMAGINARY(A10)"
$ws.Range("G10").Formula = "=IMABS(A10)"
$ws.Range("I10").Formula = "=IMARGUMENT(A10)"

$ws.Range("A11").Formula = "=COMPLEX(0,0)"
$ws.Range("C11").Formula = "=IMREAL(A11)"
$ws.Range("E11").Formula = "=IMAGINARY(A11)"
$ws.Range("G11").Formula = "=IMABS(A11)"
$ws.Range("I11").Formula = "=IMARGUMENT(A11)"

# Block 4: real power change-matrix entries (rows 13-15)
$ws.Range("A13").Value = "real power by real part"
$ws.Range("E13").Value = "real power by imaginary part"
$ws.Range("I13").Value = "real power by angle"

$ws.Range("A14").Formula = "=E2*C7+2*D2*C6-C10"
$ws.Range("A14").NumberFormat = "0.00"
$ws.Range("E14").Formula = "=H2*C7-E10"
$ws.Range("I14").Formula = "=G6*K2*G7*SIN(-N2)"

$ws.Range("A15").Formula = "=C7*D3"
$ws.Range("E15").Formula = "=-C7*G3"
$ws.Range("I15").Formula = "=C7*J3*C6*SIN(-M3)"

# Block 5: imaginary power change-matrix entries (rows 17-18)
$ws.Range("A17").Value = "imaginary power by real part"
$ws.Range("E17").Value = "imaginary power by imaginary part"
$ws.Range("I17").Value = "imaginary power by angle"

$ws.Range("A18").Formula = "=-H2*C7-2*G2*C6+E10"
$ws.Range("E18").Formula = "=E2*C7-C10"
$ws.Range("I18").Formula = "=-C6*K2*C7*COS(-N2)"

# --- 3. Final view state: new sheet active, cell L18 selected --------------
$ws.Range("L18").Select()
$ws.Activate()
